# Updated cryptos list on Sun Jun 16 14:41:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.600.92"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'3.584.59"
$ws.Range("E3").Value = "  +0.76%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'609.05"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").Value = "'146.23"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").Value = "'3.581.65"
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("D11").Value = "'7.95"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "'0.416"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "'4.193.98"
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("E14").Value = "  +0.56%  "
$ws.Range("D15").Value = "'30.04"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").Value = "'3.582.75"
$ws.Range("E16").Value = "  +0.75%  "
$ws.Range("D17").Value = "'66.653.10"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'11.48"
$ws.Range("E19").Value = "  -1.13%  "
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "'14.97"
$ws.Range("E21").Value = "  +0.26%  "
$ws.Range("D22").Value = "'432.76"
$ws.Range("E22").Value = "  +0.92%  "
$ws.Range("D23").Value = "'0.622"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").Value = "'79.06"
$ws.Range("E24").Value = "  +0.49%  "
$ws.Range("D25").Value = "'3.730.27"
$ws.Range("E25").Value = "  +0.99%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "'8.07"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +1.45%  "
$ws.Range("E31").Value = "  +0.12%  "
$ws.Range("D32").Value = "'3.580.72"
$ws.Range("D33").Value = "'25.49"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'1.71"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").Value = "'173.92"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("D43").Value = "'0.893"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +0.96%  "
$ws.Range("D45").Value = "'45.82"
$ws.Range("E45").Value = "  +0.18%  "
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").Value = "'2.52"
$ws.Range("E47").Value = "  +5.14%  "
$ws.Range("E48").Value = "  -1.39%  "
$ws.Range("D49").Value = "'25.00"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.20"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'23.62"
$ws.Range("E51").Value = "  +4.51%  "
